# edit.ps1 - apply the "finish basic structure of Mex..." commit to diagram.pptx
#
# 1) Refresh the cached "datetimeFigureOut" field text (3/2/2025 -> 3/7/2025) on
#    the slide master and all 11 slide layouts.
# 2) In slide 1's "PatternWindow" box, merge the two runs of the
#    "Initialize canvas upon window creation" bullet into a single run, and add
#    a new bullet describing the PatternMemory member.

function Set-DatePlaceholderText($shapes, $newText) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePh = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) { $isDatePh = $true }
        } catch {
        }
        if ($isDatePh) {
            # Force a full, single-run text replace: first stomp the text with
            # something unrelated (so the engine can't diff it down to a
            # shared-prefix edit), then set the real value.
            $tr = $shp.TextFrame.TextRange
            $tr.Text = "TEMP_DATE_PLACEHOLDER_VALUE"
            $tr2 = $shp.TextFrame.TextRange
            $tr2.Text = $newText
        }
    }
}

$p = $ppt.ActivePresentation

# --- 1) Date field refresh (slide master + all slide layouts) -------------
$sm = $p.Designs.Item(1).SlideMaster
Set-DatePlaceholderText $sm.Shapes "3/7/2025"
for ($L = 1; $L -le $sm.CustomLayouts.Count; $L++) {
    $lay = $sm.CustomLayouts.Item($L)
    Set-DatePlaceholderText $lay.Shapes "3/7/2025"
}

# --- 2) PatternWindow box bullet edits on slide 1 --------------------------
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(3)
$tr = $shp.TextFrame.TextRange

# Paragraph 3 currently holds two runs: "Initialize canvas upon " + "window
# creation". Collapse them into a single run (same trick as above: stomp
# then set, so the edit isn't diffed into a tiny shared-prefix patch).
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "TEMP_BULLET_PLACEHOLDER_VALUE"
$tr = $shp.TextFrame.TextRange
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "Initialize canvas upon window creation"

# Append a new bullet paragraph after it, inheriting the same bullet
# formatting (marL/indent/buFont/buChar) from paragraph 3.
$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("`rManage internal PatternMemory to store pre-loaded/dynamically-generated patterns")

# Split the new paragraph's single run into three runs so "PatternMemory"
# can carry its own run properties, matching the authored formatting.
$tr = $shp.TextFrame.TextRange
$para4 = $tr.Paragraphs(4, 1)
$full = $para4.Text
$markStart = $full.IndexOf("PatternMemory") + 1
$markLen = "PatternMemory".Length
$mid = $para4.Characters($markStart, $markLen)
$mid.Font.Size = 12
